$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.562.37"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "1.766.98"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.61%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.36"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3871"
$ws.Range("E7").Value = "  +2.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3421"
$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.19"
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.148"
$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07464"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.58"
$ws.Range("E13").Value = "  +4.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.389"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").Value = "1.773.86"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.095"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001079"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06697"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.72"
$ws.Range("E19").Value = "  -1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.49"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.463"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").Value = "27.583.84"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.20"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.85"
$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.457"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.437"
$ws.Range("E28").Value = "  -4.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.71"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.98"
$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("D31").Value = "1.975.87"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.177"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.968"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08879"
$ws.Range("E34").Value = "  +2.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.81"
$ws.Range("E35").Value = "  -3.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02444"
$ws.Range("E36").Value = "  +4.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.415"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6828"
$ws.Range("E38").Value = "  -1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06376"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2203"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.554"
$ws.Range("E41").Value = "  -6.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.246"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.424"
$ws.Range("E43").Value = "  -4.18%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.007"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.18"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6288"
$ws.Range("E46").Value = "  -3.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.856"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.02"
$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.115"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07421"
$ws.Range("E50").Value = "  +4.15%  "

$ws.Range("B51").Value = "Tezos"
$ws.Range("C51").Value = "https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.244"
$ws.Range("E51").Value = "  +3.01%  "
